$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped from
# 45172 (2023-09-03) to 45175 (2023-09-06) for every data row (2-452).
$ws.Range("C2:C452").Value = 45175
